$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40" (rule label). The edit turns it into
# the text "1" while keeping the cell's existing formatting (style) intact.
#
# A plain  $ws.Range("B11").Value = "1"  would get auto-coerced to the
# *number* 1 by Excel's type inference (since "1" looks numeric), and
# forcing text via NumberFormat/quote-prefix bumps the cell onto a new
# style record. To avoid any visible style change we snapshot the cell's
# current formatting, overwrite the value as text, then restore the
# snapshot's formatting back onto the cell.

$backup = $ws.Range("Z1")
$backup.Value = $null
$ws.Range("B11").Copy($backup)        # back up B11's current formatting

$ws.Range("B11").NumberFormat = "@"   # force text interpretation
$ws.Range("B11").Value = "1"          # write the new text value

$backup.Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats: restore original formatting
$excel.CutCopyMode = $false
$backup.Clear()                       # remove the temporary helper cell entirely
